$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a whole paragraph's contents (pPr + runs) via InsertXML.
# $range must be a Range that spans a *complete* paragraph, end-mark included
# (e.g. a Paragraph's .Range, or a union of several whole paragraphs), since
# InsertXML only *replaces* content when the target range covers full
# paragraph(s); on a sub-paragraph range it merely inserts alongside.
# ---------------------------------------------------------------------------
function Replace-ParagraphsXml($range, [string]$bodyInnerXml) {
    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
__BODY__
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
    $xml = $xml.Replace("__BODY__", $bodyInnerXml)
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "Front-end Design:" (Heading2) -> two runs: "Code" + ":"
# ---------------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("Front-end Design:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find paragraph 'Front-end Design:'"
}
$headingPara = $find1.Paragraphs(1)
$headingRange = $headingPara.Range

$headingBody = @'
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:t>Code</w:t>
            </w:r>
            <w:r>
              <w:t>:</w:t>
            </w:r>
          </w:p>
'@
Replace-ParagraphsXml $headingRange $headingBody

# ---------------------------------------------------------------------------
# 2) The four empty paragraphs right after "... Lấy từ Project Java" become:
#      - "Database: Lấy từ Project Java"  (list item)
#      - "Code: Translate code từ Project Java sang (...)" (list item)
#      - three empty paragraphs
#      - "Kế hoạch:" (Heading2)
#      - "Chuyển code theo Design Pattern" (list item)
# ---------------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("Project Java", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find anchor text 'Project Java'"
}
$anchorPara = $find2.Paragraphs(1)
$firstEmpty = $anchorPara.Next(1)
$lastEmpty = $anchorPara.Next(4)
$emptySpan = $d.Range($firstEmpty.Range.Start, $lastEmpty.Range.End)

$newBody = @'
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="0070C0"/>
              </w:rPr>
              <w:t xml:space="preserve">Database: </w:t>
            </w:r>
            <w:r>
              <w:t>L&#7845;y t&#7915; Project Java</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:color w:val="0070C0"/>
              </w:rPr>
              <w:t>Code</w:t>
            </w:r>
            <w:r>
              <w:t>: Translate code t&#7915; Project Java sang (front-end, controller, t&#7915;ng function v&#224; n&#7897;i dung c&#225;c function trong controller, code DAO,...)</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p/>
          <w:p/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>K&#7871; ho&#7841;ch:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Chuy&#7875;n code theo Design Pattern</w:t>
            </w:r>
          </w:p>
'@
Replace-ParagraphsXml $emptySpan $newBody

Write-Host "Edit complete"
